$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3091.1
$ws.Range("I40").Value = 2245.6924
$ws.Range("K40").Value = 2245.6924
$ws.Range("M40").Value = -2070.6924

$ws.Range("H45").Value = 7552.857
$ws.Range("J45").Value = 10290.6
$ws.Range("L45").Value = 30871.8
$ws.Range("N45").Value = -31255.8

$ws.Range("H64").Value = 7966.6113
$ws.Range("I64").Value = 6122.25
$ws.Range("J64").Value = 8493.571
$ws.Range("K64").Value = 6122.25
$ws.Range("L64").Value = 8493.571
$ws.Range("M64").Value = -5874.25
$ws.Range("N64").Value = -8989.571

$ws.Range("H67").Value = 7966.6113
$ws.Range("I67").Value = 6122.25
$ws.Range("J67").Value = 8493.571
$ws.Range("K67").Value = 6122.25
$ws.Range("L67").Value = 8493.571
$ws.Range("M67").Value = -5264.25
$ws.Range("N67").Value = -10209.571

$ws.Range("H76").Value = 8600.799999999999
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 8600.799999999999
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws.Range("H98").Value = 4033.394
$ws.Range("I98").Value = 1046.36
$ws.Range("J98").Value = 13367.875
$ws.Range("K98").Value = 1046.36
$ws.Range("L98").Value = 13367.875
$ws.Range("M98").Value = 451.6400000000001
$ws.Range("N98").Value = -16363.875

$ws.Range("H106").Value = 13244
$ws.Range("I106").Value = 4415.143
$ws.Range("J106").Value = 20110.889
$ws.Range("K106").Value = 4415.143
$ws.Range("L106").Value = 20110.889
$ws.Range("M106").Value = -3784.143
$ws.Range("N106").Value = -21372.889

$ws.Range("H112").Value = 2208.111
$ws.Range("I112").Value = 1512
$ws.Range("K112").Value = 4536
$ws.Range("M112").Value = -3428

$ws.Range("H116").Value = 8036.516
$ws.Range("I116").Value = 9404.333000000001
$ws.Range("K116").Value = 9404.333000000001
$ws.Range("M116").Value = -5962.333000000001

$ws.Range("H122").Value = 4033.394
$ws.Range("I122").Value = 1046.36
$ws.Range("J122").Value = 13367.875
$ws.Range("K122").Value = 3139.08
$ws.Range("L122").Value = 40103.625
$ws.Range("M122").Value = -689.0799999999999
$ws.Range("N122").Value = -45003.625

$ws.Range("H132").Value = 1519.4
$ws.Range("I132").Value = 1519.4
$ws.Range("K132").Value = 4558.200000000001
$ws.Range("M132").Value = -2028.200000000001

$ws.Range("H135").Value = 768.6316
$ws.Range("I135").Value = 660.5625
$ws.Range("K135").Value = 5945.0625
$ws.Range("M135").Value = -3410.0625

$ws.Range("H137").Value = 21279322
$ws.Range("I137").Value = 40001816
$ws.Range("J137").Value = 3762.818
$ws.Range("K137").Value = 120005448
$ws.Range("L137").Value = 11288.454
$ws.Range("M137").Value = -120002898
$ws.Range("N137").Value = -16388.454

$ws.Range("H138").Value = 4537.1787
$ws.Range("I138").Value = 2623.5
$ws.Range("J138").Value = 5302.65
$ws.Range("K138").Value = 7870.5
$ws.Range("L138").Value = 15907.95
$ws.Range("M138").Value = -2730.5
$ws.Range("N138").Value = -26187.95

$ws.Range("H141").Value = 944.875
$ws.Range("I141").Value = 944.875
$ws.Range("K141").Value = 2834.625
$ws.Range("M141").Value = 2345.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 348.66666
$ws.Range("I5").Value = 319.6
$ws.Range("J5").Value = 494
$ws.Range("K5").Value = 319.6
$ws.Range("L5").Value = 494
$ws.Range("M5").Value = -207.6
$ws.Range("N5").Value = -718

$ws.Range("H32").Value = 2162.3384
$ws.Range("I32").Value = 1987.5873
$ws.Range("K32").Value = 1987.5873
$ws.Range("M32").Value = -1700.5873

$ws.Range("H111").Value = 73683.2
$ws.Range("J111").Value = 73683.2
$ws.Range("L111").Value = 73683.2
$ws.Range("N111").Value = -81863.2

$ws.Range("H122").Value = 76925810
$ws.Range("I122").Value = 2956.5833
$ws.Range("K122").Value = 8869.749899999999
$ws.Range("M122").Value = -6419.749899999999

$ws.Range("H132").Value = 9887.375
$ws.Range("I132").Value = 5077.2856
$ws.Range("J132").Value = 16621.5
$ws.Range("K132").Value = 15231.8568
$ws.Range("L132").Value = 49864.5
$ws.Range("M132").Value = -12701.8568
$ws.Range("N132").Value = -54924.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 348.66666
$ws.Range("I4").Value = 319.6
$ws.Range("J4").Value = 494
$ws.Range("K4").Value = 319.6
$ws.Range("L4").Value = 494
$ws.Range("M4").Value = -204.6
$ws.Range("N4").Value = -724

$ws.Range("H20").Value = 2980.6897
$ws.Range("I20").Value = 2022.619
$ws.Range("K20").Value = 2022.619
$ws.Range("M20").Value = -1775.619

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3896.6
$ws.Range("I58").Value = 1982.0588
$ws.Range("J58").Value = 7965
$ws.Range("K58").Value = 1982.0588
$ws.Range("L58").Value = 7965
$ws.Range("M58").Value = -1779.0588
$ws.Range("N58").Value = -8371

$ws.Range("H82").Value = 42300
$ws.Range("I82").Value = 35000
$ws.Range("K82").Value = 35000
$ws.Range("M82").Value = -34639

$ws.Range("H85").Value = 42300
$ws.Range("I85").Value = 35000
$ws.Range("K85").Value = 35000
$ws.Range("M85").Value = -33752

$ws.Range("H107").Value = 910.6667
$ws.Range("I107").Value = 830.2
$ws.Range("J107").Value = 1011.25
$ws.Range("K107").Value = 830.2
$ws.Range("L107").Value = 1011.25
$ws.Range("M107").Value = 1089.8
$ws.Range("N107").Value = -4851.25

$ws.Range("H132").Value = 4247.56
$ws.Range("I132").Value = 2019.4651
$ws.Range("J132").Value = 17934.428
$ws.Range("K132").Value = 6058.3953
$ws.Range("L132").Value = 53803.284
$ws.Range("M132").Value = -3528.3953
$ws.Range("N132").Value = -58863.284

$ws.Range("H134").Value = 4843.3096
$ws.Range("I134").Value = 2445.0967
$ws.Range("K134").Value = 7335.2901
$ws.Range("M134").Value = -4800.2901

$ws.Range("H136").Value = 3896.6
$ws.Range("I136").Value = 1982.0588
$ws.Range("J136").Value = 7965
$ws.Range("K136").Value = 5946.1764
$ws.Range("L136").Value = 23895
$ws.Range("M136").Value = -3396.1764
$ws.Range("N136").Value = -28995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 627069.25
$ws.Range("I3").Value = 1250487.5
$ws.Range("J3").Value = 3651
$ws.Range("K3").Value = 1250487.5
$ws.Range("L3").Value = 3651
$ws.Range("M3").Value = -1250371.5
$ws.Range("N3").Value = -3883

$ws.Range("H10").Value = 13751.25
$ws.Range("I10").Value = 13751.25
$ws.Range("K10").Value = 13751.25
$ws.Range("M10").Value = -13582.25

$ws.Range("H126").Value = 2504.2856
$ws.Range("I126").Value = 1093.8182
$ws.Range("J126").Value = 7676
$ws.Range("K126").Value = 3281.4546
$ws.Range("L126").Value = 23028
$ws.Range("M126").Value = -811.4546
$ws.Range("N126").Value = -27968

$ws.Range("H132").Value = 339784.9
$ws.Range("I132").Value = 388814.66
$ws.Range("J132").Value = 6382.6
$ws.Range("K132").Value = 1166443.98
$ws.Range("L132").Value = 19147.8
$ws.Range("M132").Value = -1163913.98
$ws.Range("N132").Value = -24207.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8164.923
$ws.Range("I40").Value = 7261.5835
$ws.Range("K40").Value = 7261.5835
$ws.Range("M40").Value = -7125.5835

$ws.Range("H61").Value = 6809.364
$ws.Range("I61").Value = 3570.7144
$ws.Range("K61").Value = 3570.7144
$ws.Range("M61").Value = -3368.7144

$ws.Range("H113").Value = 6809.364
$ws.Range("I113").Value = 3570.7144
$ws.Range("K113").Value = 3570.7144
$ws.Range("M113").Value = -1400.7144

$ws.Range("H132").Value = 2971.5
$ws.Range("I132").Value = 1571.2222
$ws.Range("K132").Value = 4713.6666
$ws.Range("M132").Value = -2183.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H132").Value = 5430.1626
$ws.Range("I132").Value = 2226.3547
$ws.Range("K132").Value = 6679.0641
$ws.Range("M132").Value = -4149.0641

$ws.Range("H136").Value = 4186.2173
$ws.Range("I136").Value = 2163.9
$ws.Range("J136").Value = 17668.334
$ws.Range("K136").Value = 6491.700000000001
$ws.Range("L136").Value = 53005.00199999999
$ws.Range("M136").Value = -3941.700000000001
$ws.Range("N136").Value = -58105.00199999999
